# The commit re-deploys the rendered document; the only substantive change
# recorded in the diff is that every bookmark in the doc got a brand-new
# internal w:id (the bookmark *names* and everything else are untouched).
#
# Word's object model does not expose the raw bookmarkStart/bookmarkEnd
# "w:id" attribute as a settable property (Bookmark has no .ID / .Id member
# - that's purely an internal OOXML plumbing detail). The supported way to
# force a bookmark to be assigned a fresh internal id while preserving its
# name and location is to delete it and re-Add it over the same Range.
$d = $word.ActiveDocument

$bookmarkNames = @("pest_table", "ind_plots", "dv_vs_pred_ipred", "prm_vs_iteration")

foreach ($name in $bookmarkNames) {
    if ($d.Bookmarks.Exists($name)) {
        $bm = $d.Bookmarks($name)
        $rng = $bm.Range
        $bm.Delete()
        $d.Bookmarks.Add($name, $rng)
    }
}
